$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.46"
$ws.Range("D3").Value = "'21.73"
$ws.Range("D4").Value = "'5.482"
$ws.Range("D6").Value = "'3.373"
$ws.Range("D7").Value = "'0.8053"
$ws.Range("D8").Value = "'1.039"
$ws.Range("D9").Value = "'0.1501"
$ws.Range("D10").Value = "'0.07421"
$ws.Range("D11").Value = "'0.03175"
$ws.Range("D12").Value = "'0.03016"
$ws.Range("D13").Value = "'0.09295"
$ws.Range("D14").Value = "'3.436"
$ws.Range("D16").Value = "'0.04710"
$ws.Range("D17").Value = "'0.0005860"
$ws.Range("D18").Value = "'0.006349"
$ws.Range("D19").Value = "'0.005052"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("D20").Value = "'0.001042"
$ws.Range("D22").Value = "'0.0003200"
$ws.Range("D23").Value = "'3.766"
$ws.Range("D25").Value = "'2.151"
$ws.Range("D26").Value = "'0.3279"
$ws.Range("D40").Value = "'0.04121"
$ws.Range("D41").Value = "'0.006961"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D44").Value = "'0.009130"
$ws.Range("D45").Value = "'0.00005841"
$ws.Range("D47").Value = "'0.0005500"
$ws.Range("D48").Value = "'0.6824"
$ws.Range("D49").Value = "'0.009250"
